$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $text) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "260.80"
Set-TextValue "E2" "1.68%"
Set-TextValue "D3" "27.42"
Set-TextValue "E3" "2.40%"
Set-TextValue "D4" "4.681"
Set-TextValue "E4" "0.21%"
Set-TextValue "D5" "0.06092"
Set-TextValue "E5" "2.52%"
Set-TextValue "D6" "6.665"
Set-TextValue "E6" "0.78%"
Set-TextValue "D7" "0.8488"
Set-TextValue "E7" "-0.23%"
Set-TextValue "D8" "0.9239"
Set-TextValue "E8" "1.20%"
Set-TextValue "D9" "0.1407"
Set-TextValue "E9" "2.16%"
Set-TextValue "D10" "0.04865"
Set-TextValue "E10" "6.65%"
Set-TextValue "D11" "0.07098"
Set-TextValue "E11" "1.39%"
Set-TextValue "D12" "0.03077"
Set-TextValue "E12" "0.71%"
Set-TextValue "D13" "0.09060"
Set-TextValue "E13" "-0.31%"
Set-TextValue "D14" "0.001536"
Set-TextValue "E14" "1.03%"
Set-TextValue "D15" "0.0006094"
Set-TextValue "E15" "-94.04%"
Set-TextValue "D16" "0.006115"
Set-TextValue "E16" "1.58%"
Set-TextValue "D17" "3.450"
Set-TextValue "E17" "-0.54%"
Set-TextValue "D18" "3.148"
Set-TextValue "E18" "-0.36%"
Set-TextValue "E19" "-0.65%"
Set-TextValue "E20" "2.72%"
Set-TextValue "D21" "0.1310"
Set-TextValue "E21" "1.20%"
Set-TextValue "D22" "4.087"
Set-TextValue "E22" "5.45%"
Set-TextValue "E23" "-0.06%"
Set-TextValue "E24" "0.63%"
Set-TextValue "D25" "0.003800"
Set-TextValue "E25" "-20.24%"
Set-TextValue "E26" "0.06%"
Set-TextValue "E27" "3.42%"
Set-TextValue "D40" "0.03856"
Set-TextValue "E40" "2.27%"
Set-TextValue "D41" "0.1112"
Set-TextValue "E41" "1.60%"
Set-TextValue "D42" "0.004071"
Set-TextValue "E42" "-34.33%"
Set-TextValue "D43" "0.01623"
Set-TextValue "E43" "11.87%"
Set-TextValue "E44" "-4.46%"
Set-TextValue "D45" "0.00005143"
Set-TextValue "E45" "-2.15%"
Set-TextValue "E46" "0.06%"
Set-TextValue "E47" "-43.24%"
Set-TextValue "E48" "26.67%"
Set-TextValue "E49" "0.06%"
Set-TextValue "D50" "0.0002001"
Set-TextValue "E50" "0.06%"
